$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 22819.727
$ws.Range("J21").Value = 20142.857
$ws.Range("L21").Value = 20142.857
$ws.Range("N21").Value = -21078.857
$ws.Range("H23").Value = 22819.727
$ws.Range("J23").Value = 20142.857
$ws.Range("L23").Value = 20142.857
$ws.Range("N23").Value = -20610.857
$ws.Range("H29").Value = 4000
$ws.Range("J29").Value = 4000
$ws.Range("L29").Value = 12000
$ws.Range("N29").Value = -12562
$ws.Range("H38").Value = 682.8182
$ws.Range("I38").Value = 151.375
$ws.Range("J38").Value = 2100
$ws.Range("K38").Value = 454.125
$ws.Range("L38").Value = 6300
$ws.Range("M38").Value = -82.125
$ws.Range("N38").Value = -7044
$ws.Range("H41").Value = 10526589
$ws.Range("I41").Value = 333.6
$ws.Range("J41").Value = 22222428
$ws.Range("K41").Value = 333.6
$ws.Range("L41").Value = 22222428
$ws.Range("M41").Value = 106.4
$ws.Range("N41").Value = -22223308
$ws.Range("H58").Value = 68933.336
$ws.Range("I58").Value = 1815
$ws.Range("J58").Value = 79259.234
$ws.Range("K58").Value = 5445
$ws.Range("L58").Value = 237777.702
$ws.Range("N58").Value = -238077.702
$ws.Range("M58").Value = -5295
$ws.Range("H86").Value = 2016.6666
$ws.Range("I86").Value = 1950
$ws.Range("J86").Value = 2050
$ws.Range("K86").Value = 1950
$ws.Range("L86").Value = 2050
$ws.Range("M86").Value = -827
$ws.Range("N86").Value = -4296
$ws.Range("H87").Value = 29230.965
$ws.Range("J87").Value = 29230.965
$ws.Range("L87").Value = 29230.965
$ws.Range("N87").Value = -31726.965
$ws.Range("H89").Value = 2016.6666
$ws.Range("I89").Value = 1950
$ws.Range("J89").Value = 2050
$ws.Range("K89").Value = 9750
$ws.Range("L89").Value = 10250
$ws.Range("M89").Value = -4134
$ws.Range("N89").Value = -21482
$ws.Range("H90").Value = 29230.965
$ws.Range("J90").Value = 29230.965
$ws.Range("L90").Value = 87692.895
$ws.Range("N90").Value = -100172.895
$ws.Range("H135").Value = 398.67648
$ws.Range("I135").Value = 289.54544
$ws.Range("K135").Value = 2605.90896
$ws.Range("M135").Value = -70.90895999999975
$ws.Range("H138").Value = 4165.573
$ws.Range("I138").Value = 1514.8485
$ws.Range("J138").Value = 6248.2856
$ws.Range("K138").Value = 4544.5455
$ws.Range("L138").Value = 18744.8568
$ws.Range("M138").Value = 595.4544999999998
$ws.Range("N138").Value = -29024.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4653.1
$ws.Range("I32").Value = 4267.931
$ws.Range("J32").Value = 7230.769
$ws.Range("K32").Value = 4267.931
$ws.Range("L32").Value = 7230.769
$ws.Range("M32").Value = -3980.931
$ws.Range("N32").Value = -7804.769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3722.6206
$ws.Range("I62").Value = 2363.5557
$ws.Range("K62").Value = 2363.5557
$ws.Range("M62").Value = -1739.5557
$ws.Range("H65").Value = 3722.6206
$ws.Range("I65").Value = 2363.5557
$ws.Range("K65").Value = 11817.7785
$ws.Range("M65").Value = -8697.7785
$ws.Range("H86").Value = 2917.1177
$ws.Range("I86").Value = 1739.3
$ws.Range("K86").Value = 1739.3
$ws.Range("M86").Value = -616.3
$ws.Range("H89").Value = 2917.1177
$ws.Range("I89").Value = 1739.3
$ws.Range("K89").Value = 8696.5
$ws.Range("M89").Value = -3080.5
$ws.Range("H122").Value = 2925.238
$ws.Range("I122").Value = 2395.647
$ws.Range("J122").Value = 5176
$ws.Range("K122").Value = 7186.941
$ws.Range("L122").Value = 15528
$ws.Range("M122").Value = -4736.941
$ws.Range("N122").Value = -20428

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 34020.062
$ws.Range("I129").Value = 3488.75
$ws.Range("J129").Value = 64551.375
$ws.Range("K129").Value = 10466.25
$ws.Range("L129").Value = 193654.125
$ws.Range("M129").Value = -5466.25
$ws.Range("N129").Value = -203654.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3918.3914
$ws.Range("I80").Value = 3408.7856
$ws.Range("J80").Value = 4711.1113
$ws.Range("K80").Value = 3408.7856
$ws.Range("L80").Value = 4711.1113
$ws.Range("M80").Value = -2410.7856
$ws.Range("N80").Value = -6707.1113
$ws.Range("H83").Value = 3918.3914
$ws.Range("I83").Value = 3408.7856
$ws.Range("J83").Value = 4711.1113
$ws.Range("K83").Value = 17043.928
$ws.Range("L83").Value = 23555.5565
$ws.Range("M83").Value = -12051.928
$ws.Range("N83").Value = -33539.5565
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2360.8948
$ws.Range("I7").Value = 1960.8572
$ws.Range("J7").Value = 3481
$ws.Range("K7").Value = 1960.8572
$ws.Range("L7").Value = 3481
$ws.Range("M7").Value = -1848.8572
$ws.Range("N7").Value = -3705
$ws.Range("H22").Value = 923.2273
$ws.Range("I22").Value = 432.6
$ws.Range("K22").Value = 432.6
$ws.Range("M22").Value = -137.6
$ws.Range("H27").Value = 923.2273
$ws.Range("I27").Value = 432.6
$ws.Range("K27").Value = 432.6
$ws.Range("M27").Value = -325.6
$ws.Range("H46").Value = 1814.2858
$ws.Range("I46").Value = 300
$ws.Range("J46").Value = 2227.2727
$ws.Range("K46").Value = 300
$ws.Range("L46").Value = 2227.2727
$ws.Range("M46").Value = -112
$ws.Range("N46").Value = -2603.2727
$ws.Range("H68").Value = 2552.7273
$ws.Range("I68").Value = 1035
$ws.Range("J68").Value = 6600
$ws.Range("K68").Value = 1035
$ws.Range("L68").Value = 6600
$ws.Range("M68").Value = -286
$ws.Range("N68").Value = -8098
$ws.Range("H71").Value = 2552.7273
$ws.Range("I71").Value = 1035
$ws.Range("J71").Value = 6600
$ws.Range("K71").Value = 5175
$ws.Range("L71").Value = 33000
$ws.Range("M71").Value = -1431
$ws.Range("N71").Value = -40488
$ws.Range("H82").Value = 2887.0833
$ws.Range("I82").Value = 2298
$ws.Range("J82").Value = 4654.3335
$ws.Range("K82").Value = 2298
$ws.Range("L82").Value = 4654.3335
$ws.Range("M82").Value = -1937
$ws.Range("N82").Value = -5376.3335
$ws.Range("H85").Value = 2887.0833
$ws.Range("I85").Value = 2298
$ws.Range("J85").Value = 4654.3335
$ws.Range("K85").Value = 2298
$ws.Range("L85").Value = 4654.3335
$ws.Range("M85").Value = -1050
$ws.Range("N85").Value = -7150.3335
$ws.Range("H126").Value = 2360.8948
$ws.Range("I126").Value = 1960.8572
$ws.Range("J126").Value = 3481
$ws.Range("K126").Value = 5882.571599999999
$ws.Range("L126").Value = 10443
$ws.Range("M126").Value = -3412.571599999999
$ws.Range("N126").Value = -15383
$ws.Range("H132").Value = 2553.3257
$ws.Range("I132").Value = 1837.9615
$ws.Range("J132").Value = 3647.4119
$ws.Range("K132").Value = 5513.8845
$ws.Range("L132").Value = 10942.2357
$ws.Range("M132").Value = -2983.8845
$ws.Range("N132").Value = -16002.2357
$ws.Range("H136").Value = 4747.25
$ws.Range("I136").Value = 4261.2
$ws.Range("J136").Value = 5308.077
$ws.Range("K136").Value = 12783.6
$ws.Range("L136").Value = 15924.231
$ws.Range("M136").Value = -10233.6
$ws.Range("N136").Value = -21024.231

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 13211.857
$ws.Range("I132").Value = 3802.7742
$ws.Range("K132").Value = 11408.3226
$ws.Range("M132").Value = -8878.3226
$ws.Range("H138").Value = 29553.182
$ws.Range("J138").Value = 29553.182
$ws.Range("L138").Value = 29553.182
$ws.Range("N138").Value = -39833.182
